$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 is the "PCB1" BOM line. It previously had no pricing data in
# columns I-L (Unit Price @ Qty 1 / Price @ Qty 1 / Unit Price @ Qty 100 /
# Price @ Qty 100). Fill those cells in, following the same formula
# patterns used by every other row in the table.

# I10: Unit Price @ Qty 1 = 18/5
$ws.Range("I10").Formula = "=18/5"
# J10: Price @ Qty 1 = Qty * Unit Price @ Qty 1 (same pattern as J3:J17)
$ws.Range("J10").Formula = "=A10*I10"
# K10: Unit Price @ Qty 100 = 49.5/100
$ws.Range("K10").Formula = "=49.5/100"
# L10: Price @ Qty 100 = Qty * Unit Price @ Qty 100 (same pattern as L3:L17)
$ws.Range("L10").Formula = "=A10*K10"

# --- Formatting to match the rest of the pricing columns ---
# I10 / K10 (unit price cells) use the 5-decimal currency format, right
# aligned, themed Arial font - matching the other Unit Price cells (I3, K3, ...).
$ws.Range("I10").NumberFormat = '"$"#,##0.00000'
$ws.Range("I10").Font.ThemeColor = 1
$ws.Range("I10").HorizontalAlignment = -4152

$ws.Range("K10").NumberFormat = '"$"#,##0.00000'
$ws.Range("K10").Font.ThemeColor = 1
$ws.Range("K10").HorizontalAlignment = -4152

# J10 / L10 (computed price cells) use the 2-decimal currency format with
# general alignment - matching the other Price cells (J3, L3, ...).
$ws.Range("J10").NumberFormat = '"$"#,##0.00'
$ws.Range("J10").Font.ThemeColor = 1
$ws.Range("J10").HorizontalAlignment = 1

$ws.Range("L10").NumberFormat = '"$"#,##0.00'
$ws.Range("L10").Font.ThemeColor = 1
$ws.Range("L10").HorizontalAlignment = 1
